$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M31").Value = 5375.76
$ws1.Range("D34").Value = 3944.64
$ws1.Range("L47").Value = 447.79

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F31").Value = 8178.36
$ws2.Range("F34").Value = 7733.62
$ws2.Range("F47").Value = 879.99
$ws2.Range("F60").Value = 46982.24000000001

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3: 240X80 PORCELANATO
$ws3.Range("D3").Value = 7826.98
$ws3.Range("E3").Value = 12560.4974217135
$ws3.Range("F3").Value = 0.3839111547789598

# Row 11: PIEDRA SINTERIZADA
$ws3.Range("D11").Value = 12734.82
$ws3.Range("E11").Value = 6838.240249249699
$ws3.Range("F11").Value = 0.650629990294347

# Row 12: PORCELANATO
$ws3.Range("D12").Value = 16900.25
$ws3.Range("E12").Value = 31723.81
$ws3.Range("F12").Value = 0.3475697010903656

# Row 14: TOTAL
$ws3.Range("D14").Value = 50996.18
$ws3.Range("E14").Value = 48901.81284188786
$ws3.Range("F14").Value = 0.5104825287202066
